# Scheduled runner update: refresh market-price / profit figures on the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 748.9
$ws.Range("J17").Value = 610
$ws.Range("L17").Value = 1830
$ws.Range("N17").Value = -2166
$ws.Range("H53").Value = 354
$ws.Range("I53").Value = 225
$ws.Range("K53").Value = 225
$ws.Range("M53").Value = 412
$ws.Range("H62").Value = 9469.959999999999
$ws.Range("I62").Value = 8690.5
$ws.Range("K62").Value = 8690.5
$ws.Range("M62").Value = -8066.5
$ws.Range("H65").Value = 9469.959999999999
$ws.Range("I65").Value = 8690.5
$ws.Range("K65").Value = 43452.5
$ws.Range("M65").Value = -40332.5
$ws.Range("H113").Value = 57654.3
$ws.Range("I113").Value = 88111.5
$ws.Range("J113").Value = 11968.5
$ws.Range("K113").Value = 88111.5
$ws.Range("L113").Value = 11968.5
$ws.Range("M113").Value = -84857.5
$ws.Range("N113").Value = -18476.5
$ws.Range("H135").Value = 1734
$ws.Range("I135").Value = 1134.875
$ws.Range("K135").Value = 10213.875
$ws.Range("M135").Value = -7678.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4430.9
$ws.Range("I32").Value = 3591.7812
$ws.Range("K32").Value = 3591.7812
$ws.Range("M32").Value = -3304.7812
$ws.Range("H45").Value = 9796
$ws.Range("J45").Value = 2979
$ws.Range("L45").Value = 2979
$ws.Range("N45").Value = -3733
$ws.Range("H55").Value = 12594.5
$ws.Range("I55").Value = 5190
$ws.Range("K55").Value = 5190
$ws.Range("M55").Value = -4875
$ws.Range("H112").Value = 25359.834
$ws.Range("I112").Value = 25000
$ws.Range("K112").Value = 25000
$ws.Range("M112").Value = -23523

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2248.25
$ws.Range("J134").Value = 4700
$ws.Range("L134").Value = 14100
$ws.Range("N134").Value = -19170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52222.24
$ws.Range("I31").Value = 64149
$ws.Range("J31").Value = 14056.6
$ws.Range("K31").Value = 64149
$ws.Range("L31").Value = 14056.6
$ws.Range("M31").Value = -63854
$ws.Range("N31").Value = -14646.6
$ws.Range("H34").Value = 52222.24
$ws.Range("I34").Value = 64149
$ws.Range("J34").Value = 14056.6
$ws.Range("K34").Value = 64149
$ws.Range("L34").Value = 14056.6
$ws.Range("M34").Value = -63947
$ws.Range("N34").Value = -14460.6
$ws.Range("H58").Value = 2950.3076
$ws.Range("I58").Value = 2935.25
$ws.Range("J58").Value = 2974.4
$ws.Range("K58").Value = 2935.25
$ws.Range("L58").Value = 2974.4
$ws.Range("M58").Value = -2732.25
$ws.Range("N58").Value = -3380.4
$ws.Range("H107").Value = 1056.2368
$ws.Range("I107").Value = 1370.6875
$ws.Range("J107").Value = 827.5454999999999
$ws.Range("K107").Value = 1370.6875
$ws.Range("L107").Value = 827.5454999999999
$ws.Range("M107").Value = 549.3125
$ws.Range("N107").Value = -4667.5455
$ws.Range("H134").Value = 10306.658
$ws.Range("I134").Value = 7274.946
$ws.Range("K134").Value = 21824.838
$ws.Range("M134").Value = -19289.838
$ws.Range("H136").Value = 2950.3076
$ws.Range("I136").Value = 2935.25
$ws.Range("J136").Value = 2974.4
$ws.Range("K136").Value = 8805.75
$ws.Range("L136").Value = 8923.200000000001
$ws.Range("M136").Value = -6255.75
$ws.Range("N136").Value = -14023.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 933.3333
$ws.Range("I51").Value = 933.3333
$ws.Range("K51").Value = 2799.9999
$ws.Range("M51").Value = -2339.9999
$ws.Range("H129").Value = 657.5833
$ws.Range("I129").Value = 581
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 1743
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 3257
$ws.Range("N129").Value = -14500
$ws.Range("H131").Value = 101950.77
$ws.Range("I131").Value = 250479.88
$ws.Range("K131").Value = 751439.64
$ws.Range("M131").Value = -746399.64
$ws.Range("H132").Value = 1553.8462
$ws.Range("I132").Value = 1270
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 11430
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -8900
$ws.Range("N132").Value = -27560
$ws.Range("H137").Value = 1734.1428
$ws.Range("I137").Value = 1356.5
$ws.Range("K137").Value = 4069.5
$ws.Range("M137").Value = 1030.5
$ws.Range("H138").Value = 45467628
$ws.Range("I138").Value = 83344820
$ws.Range("J138").Value = 15000
$ws.Range("K138").Value = 250034460
$ws.Range("L138").Value = 45000
$ws.Range("M138").Value = -250029320
$ws.Range("N138").Value = -55280
$ws.Range("H140").Value = 3529.6365
$ws.Range("I140").Value = 3529.6365
$ws.Range("K140").Value = 10588.9095
$ws.Range("M140").Value = -5408.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 95000
$ws.Range("J88").Value = 95000
$ws.Range("L88").Value = 95000
$ws.Range("N88").Value = -95902
$ws.Range("H91").Value = 95000
$ws.Range("J91").Value = 95000
$ws.Range("L91").Value = 95000
$ws.Range("N91").Value = -98120
$ws.Range("H102").Value = 4667.8335
$ws.Range("I102").Value = 1670.6666
$ws.Range("K102").Value = 1670.6666
$ws.Range("M102").Value = -48.66660000000002
$ws.Range("H107").Value = 2953.682
$ws.Range("I107").Value = 2598.5334
$ws.Range("K107").Value = 2598.5334
$ws.Range("M107").Value = -678.5333999999998
$ws.Range("H132").Value = 837231.5
$ws.Range("I132").Value = 1114196
$ws.Range("J132").Value = 6338
$ws.Range("K132").Value = 3342588
$ws.Range("L132").Value = 19014
$ws.Range("M132").Value = -3340058
$ws.Range("N132").Value = -24074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6311.9653
$ws.Range("J7").Value = 3812.5
$ws.Range("L7").Value = 3812.5
$ws.Range("N7").Value = -4036.5
$ws.Range("H22").Value = 1365.2307
$ws.Range("I22").Value = 1339.875
$ws.Range("J22").Value = 1405.8
$ws.Range("K22").Value = 1339.875
$ws.Range("L22").Value = 1405.8
$ws.Range("M22").Value = -1044.875
$ws.Range("N22").Value = -1995.8
$ws.Range("H27").Value = 1365.2307
$ws.Range("I27").Value = 1339.875
$ws.Range("J27").Value = 1405.8
$ws.Range("K27").Value = 1339.875
$ws.Range("L27").Value = 1405.8
$ws.Range("M27").Value = -1232.875
$ws.Range("N27").Value = -1619.8
$ws.Range("H40").Value = 5077.2666
$ws.Range("I40").Value = 4262.4165
$ws.Range("K40").Value = 4262.4165
$ws.Range("M40").Value = -4126.4165
$ws.Range("H55").Value = 675
$ws.Range("I55").Value = 675
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 675
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -502
$ws.Range("N55").ClearContents()
$ws.Range("H94").Value = 40330
$ws.Range("J94").Value = 40330
$ws.Range("L94").Value = 40330
$ws.Range("N94").Value = -41682
$ws.Range("H122").Value = 5044.846
$ws.Range("I122").Value = 4358.3
$ws.Range("K122").Value = 13074.9
$ws.Range("M122").Value = -10624.9
$ws.Range("H126").Value = 6311.9653
$ws.Range("J126").Value = 3812.5
$ws.Range("L126").Value = 11437.5
$ws.Range("N126").Value = -16377.5
$ws.Range("H140").Value = 77712.5
$ws.Range("J140").Value = 77712.5
$ws.Range("L140").Value = 77712.5
$ws.Range("N140").Value = -88072.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2942.64
$ws.Range("I132").Value = 3068.087
$ws.Range("K132").Value = 9204.261
$ws.Range("M132").Value = -6674.261
